$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value. Values must remain plain text
# (matching the original inline-string cell type), so each cell is forced to
# the "@" (Text) number format before assignment and reset to the Normal style
# afterwards to avoid leaving a stray number format applied to the cell.
$updates = @{
    "D2" = "261.69"
    "E2" = "1.25%"
    "D3" = "27.39"
    "E3" = "1.32%"
    "D4" = "4.775"
    "E4" = "3.38%"
    "D5" = "0.06070"
    "E5" = "2.72%"
    "D6" = "6.713"
    "E6" = "1.05%"
    "D7" = "0.8643"
    "D8" = "0.9249"
    "E8" = "-2.34%"
    "D9" = "0.1412"
    "E9" = "0.58%"
    "D10" = "0.05059"
    "D11" = "0.07152"
    "E11" = "0.78%"
    "E12" = "-1.27%"
    "E13" = "-0.40%"
    "D14" = "0.001530"
    "E14" = "-0.37%"
    "D15" = "0.0006066"
    "E15" = "-94.21%"
    "D16" = "0.006194"
    "E16" = "1.85%"
    "D17" = "3.451"
    "E17" = "-1.43%"
    "D18" = "3.167"
    "E18" = "-0.70%"
    "D19" = "2.195"
    "E19" = "-0.41%"
    "E20" = "2.34%"
    "E21" = "0.91%"
    "D22" = "4.095"
    "E22" = "6.77%"
    "E23" = "-0.13%"
    "E24" = "-0.24%"
    "D25" = "0.003914"
    "E25" = "-8.90%"
    "E26" = "0.03%"
    "E27" = "-18.88%"
    "D40" = "0.03879"
    "E40" = "1.40%"
    "E41" = "1.13%"
    "D42" = "0.004131"
    "E42" = "-33.58%"
    "D43" = "0.01502"
    "E43" = "4.98%"
    "E44" = "0.39%"
    "D45" = "0.00005312"
    "E45" = "-1.57%"
    "E46" = "0.03%"
    "E47" = "7.01%"
    "D48" = "0.1353"
    "E48" = "-46.22%"
    "E49" = "0.03%"
    "E50" = "0.03%"
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
    $cell.Style = "Normal"
}
